# Apply the latest cryptos-list scrape: refreshed prices/volumes and two
# re-ranked coin pairs (rows 30/31 swap Fetch.AI <-> RenzoRestakedETH,
# rows 45/46 swap InjectiveProtocol <-> ONDO).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '66.528.39' }
    @{ Cell = 'E2'; Value = '  +0.62%  ' }
    @{ Cell = 'D3'; Value = '3.599.81' }
    @{ Cell = 'E3'; Value = '  +1.11%  ' }
    @{ Cell = 'E4'; Value = '  +0.00%  ' }
    @{ Cell = 'D5'; Value = '''609.41' }
    @{ Cell = 'E5'; Value = '  +0.63%  ' }
    @{ Cell = 'D6'; Value = '''149.01' }
    @{ Cell = 'E6'; Value = '  +3.05%  ' }
    @{ Cell = 'E7'; Value = '  +0.05%  ' }
    @{ Cell = 'E8'; Value = '  -0.32%  ' }
    @{ Cell = 'D9'; Value = '''8.04' }
    @{ Cell = 'E9'; Value = '  +1.65%  ' }
    @{ Cell = 'E10'; Value = '  -0.20%  ' }
    @{ Cell = 'E11'; Value = '  +0.84%  ' }
    @{ Cell = 'D12'; Value = '4.210.51' }
    @{ Cell = 'E12'; Value = '  +1.12%  ' }
    @{ Cell = 'E13'; Value = '  +0.70%  ' }
    @{ Cell = 'D14'; Value = '''29.87' }
    @{ Cell = 'E14'; Value = '  -0.36%  ' }
    @{ Cell = 'D15'; Value = '3.589.34' }
    @{ Cell = 'E15'; Value = '  +0.88%  ' }
    @{ Cell = 'D16'; Value = '66.599.94' }
    @{ Cell = 'E16'; Value = '  +0.59%  ' }
    @{ Cell = 'E17'; Value = '  +0.77%  ' }
    @{ Cell = 'E19'; Value = '  +3.21%  ' }
    @{ Cell = 'D20'; Value = '''15.13' }
    @{ Cell = 'E20'; Value = '  +1.93%  ' }
    @{ Cell = 'D21'; Value = '''427.20' }
    @{ Cell = 'E21'; Value = '  -0.42%  ' }
    @{ Cell = 'E22'; Value = '  +1.31%  ' }
    @{ Cell = 'E23'; Value = '  -0.37%  ' }
    @{ Cell = 'E24'; Value = '  -0.04%  ' }
    @{ Cell = 'D25'; Value = '''0.0000121' }
    @{ Cell = 'E25'; Value = '  +2.87%  ' }
    @{ Cell = 'E26'; Value = '  +4.72%  ' }
    @{ Cell = 'D27'; Value = '''9.45' }
    @{ Cell = 'E27'; Value = '  +4.07%  ' }
    @{ Cell = 'E28'; Value = '  +0.74%  ' }
    @{ Cell = 'E29'; Value = '  -0.04%  ' }
    @{ Cell = 'B30'; Value = 'RenzoRestakedETH' }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth' }
    @{ Cell = 'D30'; Value = '3.597.58' }
    @{ Cell = 'E30'; Value = '  +1.17%  ' }
    @{ Cell = 'B31'; Value = 'Fetch.AI' }
    @{ Cell = 'C31'; Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D31'; Value = '''1.47' }
    @{ Cell = 'E31'; Value = '  +0.70%  ' }
    @{ Cell = 'D32'; Value = '''0.157' }
    @{ Cell = 'E32'; Value = '  +3.19%  ' }
    @{ Cell = 'D33'; Value = '''25.45' }
    @{ Cell = 'E33'; Value = '  -0.41%  ' }
    @{ Cell = 'D34'; Value = '''7.86' }
    @{ Cell = 'E34'; Value = '  -0.03%  ' }
    @{ Cell = 'D36'; Value = '''5.64' }
    @{ Cell = 'E36'; Value = '  +0.75%  ' }
    @{ Cell = 'E37'; Value = '  -2.58%  ' }
    @{ Cell = 'D38'; Value = '''177.22' }
    @{ Cell = 'E38'; Value = '  +0.77%  ' }
    @{ Cell = 'E39'; Value = '  +1.09%  ' }
    @{ Cell = 'E40'; Value = '  +0.73%  ' }
    @{ Cell = 'E41'; Value = '  +0.67%  ' }
    @{ Cell = 'E42'; Value = '  -2.13%  ' }
    @{ Cell = 'E43'; Value = '  +9.91%  ' }
    @{ Cell = 'E44'; Value = '  +0.00%  ' }
    @{ Cell = 'B45'; Value = 'ONDO' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo' }
    @{ Cell = 'D45'; Value = '''1.18' }
    @{ Cell = 'E45'; Value = '  -1.55%  ' }
    @{ Cell = 'B46'; Value = 'InjectiveProtocol' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' }
    @{ Cell = 'D46'; Value = '''25.03' }
    @{ Cell = 'E46'; Value = '  -2.92%  ' }
    @{ Cell = 'D47'; Value = '''24.09' }
    @{ Cell = 'E47'; Value = '  +2.55%  ' }
    @{ Cell = 'E48'; Value = '  +1.28%  ' }
    @{ Cell = 'E49'; Value = '  +1.19%  ' }
    @{ Cell = 'D50'; Value = '2.428.94' }
    @{ Cell = 'E50'; Value = '  +5.52%  ' }
    @{ Cell = 'E51'; Value = '  -0.71%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
